$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix column definitions: col 1 should only apply to column 1 (not 1-2) ---
$ws.Columns.Item(1).ColumnWidth = 30.7109375

# --- 2. Row 10 (Objetivos:): replace the placeholder professor name with the real objective text ---
$ws.Range("B10").Value = "Apresentar aos alunos os princípios fundamentais de engenharia do meio ambiente."
$ws.Range("C10").Value = "Apresentar aos alunos os princípios fundamentais de engenharia do meio ambiente."

# --- 3. Insert a new row at 13 (this shifts the old rows 13-21 down to 14-22, carrying
#        over their formatting/row-heights automatically) to make room for the new
#        "Docentes responsáveis:" data row ---
$ws.Rows.Item(13).Insert()

# Row 13 now needs B13/C13 styled like the other data columns (copy formats from the
# existing B10/C10 cells, which already carry the correct column styles) and then the
# professor name value; row 13 keeps no label in column A.
$ws.Range("B10").Copy($ws.Range("B13"))
$ws.Range("C10").Copy($ws.Range("C13"))
$ws.Range("B13").Value = "5840671 - Francisco José Moreira Chaves"
$ws.Range("C13").Value = "5840671 - Francisco José Moreira Chaves"
$ws.Range("A13").Clear()

# --- 4. Row 14 (Programa resumido:) now carries the short syllabus text ---
$ws.Range("B14").Value = "1 - Fundamentos da Engenharia e o Meio Ambiente. 2 - O meio ambiente aquático. 3 - O meio ambiente terrestre. 4 - O meio ambiente atmosférico ."
$ws.Range("C14").Value = "1 - Fundamentos da Engenharia e o Meio Ambiente. 2 - O meio ambiente aquático. 3 - O meio ambiente terrestre. 4 - O meio ambiente atmosférico ."

# --- 5. Row 15 (Short syllabus:) no longer carries the stray "01/01/2018" value ---
$ws.Range("B15:C15").Clear()

# --- 6. Row 16 (Programa:) now carries the full syllabus text ---
$ws.Range("B16").Value = "1 - FUNDAMENTOS: A Engenharia e o Meio Ambiente; Os Ecossistemas. A crise energética. Fontes alternativas de energia. A sustentabilidade do meio ambiente. 2 - O MEIO AMBIENTE AQUÁTICO: Composição e Propriedades; Necessidade e Utilização; Requisitos de Qualidade; Poluição.3 - O MEIO AMBIENTE TERRESTRE: Composição e Propriedades; Necessidades e Utilização; Requisitos de Qualidade; Poluição.4 - O MEIO AMBIENTE ATMOSFÉRICO: Composição e Propriedades; Requisitos de Qualidade; Poluição."
$ws.Range("C16").Value = "1 - FUNDAMENTOS: A Engenharia e o Meio Ambiente; Os Ecossistemas. A crise energética. Fontes alternativas de energia. A sustentabilidade do meio ambiente. 2 - O MEIO AMBIENTE AQUÁTICO: Composição e Propriedades; Necessidade e Utilização; Requisitos de Qualidade; Poluição.3 - O MEIO AMBIENTE TERRESTRE: Composição e Propriedades; Necessidades e Utilização; Requisitos de Qualidade; Poluição.4 - O MEIO AMBIENTE ATMOSFÉRICO: Composição e Propriedades; Requisitos de Qualidade; Poluição."

# --- 7. Row 18 (Avaliação:) no longer carries the stray professor-name value ---
$ws.Range("B18:C18").Clear()

# --- 8. Row 19 (Método:) now carries the exam-method description ---
$ws.Range("B19").Value = "Duas Provas  P1  1º bimestre e P2  2º bimestre"
$ws.Range("C19").Value = "Duas Provas  P1  1º bimestre e P2  2º bimestre"

# --- 9. Row 20 (Critério:) now carries the grading formula ---
$ws.Range("B20").Value = "MF = (P1+ P2)/2"
$ws.Range("C20").Value = "MF = (P1+ P2)/2"

# --- 10. Row 21 (Norma de recuperação:) now carries the recovery-grade formula ---
$ws.Range("B21").Value = "NF = (MF + PR)/ 2 , onde PR é uma prova de recuperação"
$ws.Range("C21").Value = "NF = (MF + PR)/ 2 , onde PR é uma prova de recuperação"

# --- 11. New row 22 (Bibliografia:) with the bibliography text ---
$ws.Range("A21").Copy($ws.Range("A22"))
$ws.Range("B21").Copy($ws.Range("B22"))
$ws.Range("C21").Copy($ws.Range("C22"))
$ws.Range("A22").Value = "Bibliografia:"
$ws.Range("B22").Value = "Braga, B.P.F., M.T.,Conejo, J.G., Porto, M.F., Veras M.S., Nucci, N., Juliano, N. e Eiger, S. Introdução à Engenharia Ambiental, Makron Books, São Paulo, 1998`nSperling, M.V. Princípios do Tratamento Biológico de Águas Residuárias. Desa-UFMG, Minas Gerais, 1996.`nBRAGA, B.et al. Introdução à Engenharia Ambiental. São Paulo: Prentice Hall, 2002, 305 p.`nVON SPERLING, M. Introdução à qualidade das águas e ao tratamento de esgotos. 2. ed. Belo Horizonte: UFMG, 1996."
$ws.Range("C22").Value = "Braga, B.P.F., M.T.,Conejo, J.G., Porto, M.F., Veras M.S., Nucci, N., Juliano, N. e Eiger, S. Introdução à Engenharia Ambiental, Makron Books, São Paulo, 1998`nSperling, M.V. Princípios do Tratamento Biológico de Águas Residuárias. Desa-UFMG, Minas Gerais, 1996.`nBRAGA, B.et al. Introdução à Engenharia Ambiental. São Paulo: Prentice Hall, 2002, 305 p.`nVON SPERLING, M. Introdução à qualidade das águas e ao tratamento de esgotos. 2. ed. Belo Horizonte: UFMG, 1996."

Write-Host "edit complete"
